$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number that was updated
# from 2023-09-17 (45186) to 2023-09-19 (45188) for every data row (2-100).
$ws.Range("C2:C100").Value = 45188
